# Fruta / hortaliza, semanal
# Insert 3 new weekly rows at the top of the Kiwi / Vega Central Mapocho de
# Santiago (O'Higgins & Curicó) block, pushing the existing rows 753-798
# down to 756-801 (dimension grows from A1:T798 to A1:T801).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows right before the current row 753.
$ws.Rows("753:755").Insert()

# --- Row 753: Especial ---------------------------------------------------
$ws.Range("A753").Value = 9
$ws.Range("B753").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C753").Value = "Metropolitana"
$ws.Range("D753").Value = 45021
$ws.Range("E753").Value = 13
$ws.Range("F753").Value = "Fruta"
$ws.Range("G753").Value = 100101
$ws.Range("H753").Value = "Berries"
$ws.Range("I753").Value = 100101007
$ws.Range("J753").Value = "Kiwi"
$ws.Range("K753").Value = "Hayward"
$ws.Range("L753").Value = "Especial"
$ws.Range("M753").Value = 250
$ws.Range("N753").Value = 10000
$ws.Range("O753").Value = 10000
$ws.Range("P753").Value = 10000
$ws.Range("Q753").Value = "$/bandeja 10 kilos"
$ws.Range("R753").Value = "Provincia de Curicó"
$ws.Range("S753").Value = 1000
$ws.Range("T753").Value = 10

# --- Row 754: Primera -----------------------------------------------------
$ws.Range("A754").Value = 9
$ws.Range("B754").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C754").Value = "Metropolitana"
$ws.Range("D754").Value = 45021
$ws.Range("E754").Value = 13
$ws.Range("F754").Value = "Fruta"
$ws.Range("G754").Value = 100101
$ws.Range("H754").Value = "Berries"
$ws.Range("I754").Value = 100101007
$ws.Range("J754").Value = "Kiwi"
$ws.Range("K754").Value = "Hayward"
$ws.Range("L754").Value = "Primera"
$ws.Range("M754").Value = 220
$ws.Range("N754").Value = 8000
$ws.Range("O754").Value = 8000
$ws.Range("P754").Value = 8000
$ws.Range("Q754").Value = "$/bandeja 10 kilos"
$ws.Range("R754").Value = "Provincia de Curicó"
$ws.Range("S754").Value = 800
$ws.Range("T754").Value = 10

# --- Row 755: Segunda -------------------------------------------------------
$ws.Range("A755").Value = 9
$ws.Range("B755").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C755").Value = "Metropolitana"
$ws.Range("D755").Value = 45021
$ws.Range("E755").Value = 13
$ws.Range("F755").Value = "Fruta"
$ws.Range("G755").Value = 100101
$ws.Range("H755").Value = "Berries"
$ws.Range("I755").Value = 100101007
$ws.Range("J755").Value = "Kiwi"
$ws.Range("K755").Value = "Hayward"
$ws.Range("L755").Value = "Segunda"
$ws.Range("M755").Value = 180
$ws.Range("N755").Value = 6000
$ws.Range("O755").Value = 6000
$ws.Range("P755").Value = 6000
$ws.Range("Q755").Value = "$/bandeja 10 kilos"
$ws.Range("R755").Value = "Provincia de Curicó"
$ws.Range("S755").Value = 600
$ws.Range("T755").Value = 10
